$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Angel', ['Token Creature — Angel', 'Flying', '4/4'])"
$ws.Range("A3").Value = "('Cleric', ['Token Creature — Cleric', '{3}{W}{B}{B}, {T}, Sacrifice this creature: Return a card named Deathpact Angel from your graveyard to the battlefield.', '1/1'])"
$ws.Range("A4").Value = "('Domri Rade Emblem', ['Emblem — Domri', 'Creatures you control have double strike, trample, hexproof, and haste.'])"
$ws.Range("A5").Value = "('Frog Lizard', ['Token Creature — Frog Lizard', '3/3'])"
$ws.Range("A6").Value = "('Horror', ['Token Creature — Horror', 'Flying', '1/1'])"
$ws.Range("A7").Value = "('Rat', ['Token Creature — Rat', '1/1'])"
$ws.Range("A8").Value = "('Soldier', ['Token Creature — Soldier', '1/1'])"
$ws.Range("A9").Value = "('Spirit', ['Token Creature — Spirit', 'Flying', '1/1'])"

$ws.Range("A10:A29").EntireRow.Delete() | Out-Null
